# Applies the edit described in the commit: on the "Segmentation" sheet,
# delete row 24 of the step-diagram table (shifting rows 25-34 up to 24-33),
# and restore the small fixed "arrow" annotation (I23/J23:J26) plus fix up
# the countdown value in A23 that the row delete does not auto-update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Segmentation")

# Delete entire row 24, shifting everything below it up by one row.
$ws.Rows.Item(24).Delete()

# The row delete carries the small fixed arrow-diagram cells (J25/J26) up
# with the rest of the range, but those annotation cells are meant to stay
# anchored in place, so restore their original text.
$ws.Cells.Item(25, 10).Value = "|"
$ws.Cells.Item(26, 10).Value = "V"

# Column A holds a literal countdown (10..0) tied to each row's position;
# after removing one row from the sequence the remaining numbers are
# renumbered down by one, starting with the top row.
$ws.Cells.Item(23, 1).Value = 9

# Restore the view state: scrolled to the top, with the selection sitting
# one row below the new last data row.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A34").Select()
